# Update vm_pu.xlsx results: change B2:B25 setpoint from 1.05 to 1.02 (380 kV case)
# and refresh the resulting bus voltage magnitudes for all affected columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.060573037864563
$ws.Cells.Item(2, 4).Value = 1.067518097449828
$ws.Cells.Item(2, 5).Value = 1.073917377219075
$ws.Cells.Item(2, 6).Value = 1.080518332184557
$ws.Cells.Item(2, 9).Value = 1.055056382191564
$ws.Cells.Item(2, 10).Value = 1.065553189663523
$ws.Cells.Item(2, 11).Value = 1.070225745759889
$ws.Cells.Item(2, 12).Value = 1.07660796823292
$ws.Cells.Item(2, 13).Value = 1.083191556995883
$ws.Cells.Item(2, 14).Value = 1.067066397413504

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.061582598243144
$ws.Cells.Item(3, 4).Value = 1.068350431050944
$ws.Cells.Item(3, 5).Value = 1.074891594871721
$ws.Cells.Item(3, 6).Value = 1.081516212415311
$ws.Cells.Item(3, 9).Value = 1.055379654381123
$ws.Cells.Item(3, 10).Value = 1.066216167035869
$ws.Cells.Item(3, 11).Value = 1.070873838941376
$ws.Cells.Item(3, 12).Value = 1.077398825500713
$ws.Cells.Item(3, 13).Value = 1.084007273836979
$ws.Cells.Item(3, 14).Value = 1.067730316289764

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.062236242156173
$ws.Cells.Item(4, 4).Value = 1.068889306224261
$ws.Cells.Item(4, 5).Value = 1.075522695370386
$ws.Cells.Item(4, 6).Value = 1.082162655969556
$ws.Cells.Item(4, 9).Value = 1.055587713987035
$ws.Cells.Item(4, 10).Value = 1.066644913337262
$ws.Cells.Item(4, 11).Value = 1.071292846995463
$ws.Cells.Item(4, 12).Value = 1.077910655308213
$ws.Cells.Item(4, 13).Value = 1.084535224703967
$ws.Cells.Item(4, 14).Value = 1.068159671460083

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.062511127072667
$ws.Cells.Item(5, 4).Value = 1.069115920119744
$ws.Cells.Item(5, 5).Value = 1.075788180631375
$ws.Cells.Item(5, 6).Value = 1.082434599036772
$ws.Cells.Item(5, 9).Value = 1.05567491391933
$ws.Cells.Item(5, 10).Value = 1.066825099121293
$ws.Cells.Item(5, 11).Value = 1.071468913108022
$ws.Cells.Item(5, 12).Value = 1.078125849876473
$ws.Cells.Item(5, 13).Value = 1.084757204966335
$ws.Cells.Item(5, 14).Value = 1.068340113128594

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.062557286898229
$ws.Cells.Item(6, 4).Value = 1.069153973711902
$ws.Cells.Item(6, 5).Value = 1.075832766752699
$ws.Cells.Item(6, 6).Value = 1.082480269896314
$ws.Cells.Item(6, 9).Value = 1.055689539425115
$ws.Cells.Item(6, 10).Value = 1.066855349619114
$ws.Cells.Item(6, 11).Value = 1.071498470385545
$ws.Cells.Item(6, 12).Value = 1.078161983219027
$ws.Cells.Item(6, 13).Value = 1.084794478148274
$ws.Cells.Item(6, 14).Value = 1.068370406585592

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.062239914817074
$ws.Cells.Item(7, 4).Value = 1.068892333973351
$ws.Cells.Item(7, 5).Value = 1.075526242128209
$ws.Cells.Item(7, 6).Value = 1.082166288987125
$ws.Cells.Item(7, 9).Value = 1.055588880211641
$ws.Cells.Item(7, 10).Value = 1.066647321221805
$ws.Cells.Item(7, 11).Value = 1.071295199933852
$ws.Cells.Item(7, 12).Value = 1.077913530665769
$ws.Cells.Item(7, 13).Value = 1.084538190699339
$ws.Cells.Item(7, 14).Value = 1.068162082764098

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.060914142261487
$ws.Cells.Item(8, 4).Value = 1.067799325853396
$ws.Cells.Item(8, 5).Value = 1.074246469714405
$ws.Cells.Item(8, 6).Value = 1.080855414963273
$ws.Cells.Item(8, 9).Value = 1.055165865016801
$ws.Cells.Item(8, 10).Value = 1.065777296179729
$ws.Cells.Item(8, 11).Value = 1.070444844254737
$ws.Cells.Item(8, 12).Value = 1.076875222578449
$ws.Cells.Item(8, 13).Value = 1.083467205427748
$ws.Cells.Item(8, 14).Value = 1.067290822186669

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.058580979287228
$ws.Cells.Item(9, 4).Value = 1.065875643806631
$ws.Cells.Item(9, 5).Value = 1.071996875880441
$ws.Cells.Item(9, 6).Value = 1.078551260605629
$ws.Cells.Item(9, 9).Value = 1.054411905038878
$ws.Cells.Item(9, 10).Value = 1.064242357922203
$ws.Cells.Item(9, 11).Value = 1.068943747393124
$ws.Cells.Item(9, 12).Value = 1.07504633179618
$ws.Cells.Item(9, 13).Value = 1.081581004989244
$ws.Cells.Item(9, 14).Value = 1.065753704140757

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.057027596973701
$ws.Cells.Item(10, 4).Value = 1.064594817352939
$ws.Cells.Item(10, 5).Value = 1.070500917453766
$ws.Cells.Item(10, 6).Value = 1.077019094262241
$ws.Cells.Item(10, 9).Value = 1.053903537785736
$ws.Cells.Item(10, 10).Value = 1.063217864188861
$ws.Cells.Item(10, 11).Value = 1.067941267449961
$ws.Cells.Item(10, 12).Value = 1.073827613833101
$ws.Cells.Item(10, 13).Value = 1.080324265573983
$ws.Cells.Item(10, 14).Value = 1.064727755508803

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.056355459644153
$ws.Cells.Item(11, 4).Value = 1.064040603030793
$ws.Cells.Item(11, 5).Value = 1.06985405614781
$ws.Cells.Item(11, 6).Value = 1.076356594129239
$ws.Cells.Item(11, 9).Value = 1.053682055374132
$ws.Cells.Item(11, 10).Value = 1.062773969605006
$ws.Cells.Item(11, 11).Value = 1.067506775992867
$ws.Cells.Item(11, 12).Value = 1.07330003313362
$ws.Cells.Item(11, 13).Value = 1.079780265520282
$ws.Cells.Item(11, 14).Value = 1.064283230543723

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.056105871537194
$ws.Cells.Item(12, 4).Value = 1.063834802920896
$ws.Cells.Item(12, 5).Value = 1.069613918872909
$ws.Cells.Item(12, 6).Value = 1.076110653844967
$ws.Cells.Item(12, 9).Value = 1.053599583541247
$ws.Cells.Item(12, 10).Value = 1.062609045567798
$ws.Cells.Item(12, 11).Value = 1.067345325372418
$ws.Cells.Item(12, 12).Value = 1.073104086699312
$ws.Cells.Item(12, 13).Value = 1.079578226883224
$ws.Cells.Item(12, 14).Value = 1.064118072295465

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.056159405695441
$ws.Cells.Item(13, 4).Value = 1.063878945024492
$ws.Cells.Item(13, 5).Value = 1.06966542297446
$ws.Cells.Item(13, 6).Value = 1.076163402442246
$ws.Cells.Item(13, 9).Value = 1.053617283228547
$ws.Cells.Item(13, 10).Value = 1.062644424217431
$ws.Cells.Item(13, 11).Value = 1.067379959835884
$ws.Cells.Item(13, 12).Value = 1.073146116935011
$ws.Cells.Item(13, 13).Value = 1.079621563614357
$ws.Cells.Item(13, 14).Value = 1.064153501186839

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.056334827109193
$ws.Cells.Item(14, 4).Value = 1.064023590305242
$ws.Cells.Item(14, 5).Value = 1.069834203550425
$ws.Cells.Item(14, 6).Value = 1.076336261731838
$ws.Cells.Item(14, 9).Value = 1.053675242374178
$ws.Cells.Item(14, 10).Value = 1.062760337782662
$ws.Cells.Item(14, 11).Value = 1.067493431678101
$ws.Cells.Item(14, 12).Value = 1.073283835717621
$ws.Cells.Item(14, 13).Value = 1.079763564384485
$ws.Cells.Item(14, 14).Value = 1.064269579362628

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.056442919799477
$ws.Cells.Item(15, 4).Value = 1.064112718982333
$ws.Cells.Item(15, 5).Value = 1.069938212855099
$ws.Cells.Item(15, 6).Value = 1.076442784855326
$ws.Cells.Item(15, 9).Value = 1.053710925973519
$ws.Cells.Item(15, 10).Value = 1.062831750430155
$ws.Cells.Item(15, 11).Value = 1.067563337340617
$ws.Cells.Item(15, 12).Value = 1.073368691544934
$ws.Cells.Item(15, 13).Value = 1.079851059371469
$ws.Cells.Item(15, 14).Value = 1.064341093424272

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.057072214817931
$ws.Cells.Item(16, 4).Value = 1.064631607065513
$ws.Cells.Item(16, 5).Value = 1.070543866509607
$ws.Cells.Item(16, 6).Value = 1.077063082048349
$ws.Cells.Item(16, 9).Value = 1.053918208301536
$ws.Cells.Item(16, 10).Value = 1.06324731811502
$ws.Cells.Item(16, 11).Value = 1.067970094608693
$ws.Cells.Item(16, 12).Value = 1.073862630468873
$ws.Cells.Item(16, 13).Value = 1.080360372842412
$ws.Cells.Item(16, 14).Value = 1.064757251262915

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.057467085737419
$ws.Cells.Item(17, 4).Value = 1.064957197413626
$ws.Cells.Item(17, 5).Value = 1.070924018192808
$ws.Cells.Item(17, 6).Value = 1.077452429811005
$ws.Cells.Item(17, 9).Value = 1.054047868208956
$ws.Cells.Item(17, 10).Value = 1.063507917621736
$ws.Cells.Item(17, 11).Value = 1.068225133079193
$ws.Cells.Item(17, 12).Value = 1.074172501172912
$ws.Cells.Item(17, 13).Value = 1.08067989945022
$ws.Cells.Item(17, 14).Value = 1.065018220850822

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.057697454270839
$ws.Cells.Item(18, 4).Value = 1.065147146517236
$ws.Cells.Item(18, 5).Value = 1.071145840881533
$ws.Cells.Item(18, 6).Value = 1.077679620221188
$ws.Cells.Item(18, 9).Value = 1.054123365779453
$ws.Cells.Item(18, 10).Value = 1.063659893655628
$ws.Cells.Item(18, 11).Value = 1.068373853005375
$ws.Cells.Item(18, 12).Value = 1.074353256224744
$ws.Cells.Item(18, 13).Value = 1.080866290888062
$ws.Cells.Item(18, 14).Value = 1.065170412708114

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.057776011948907
$ws.Cells.Item(19, 4).Value = 1.065211920631361
$ws.Cells.Item(19, 5).Value = 1.071221491357424
$ws.Cells.Item(19, 6).Value = 1.077757101595923
$ws.Cells.Item(19, 9).Value = 1.054149086294008
$ws.Cells.Item(19, 10).Value = 1.063711708917158
$ws.Cells.Item(19, 11).Value = 1.068424555904981
$ws.Cells.Item(19, 12).Value = 1.074414891162454
$ws.Cells.Item(19, 13).Value = 1.080929848409274
$ws.Cells.Item(19, 14).Value = 1.065222301553259

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.057424714959773
$ws.Cells.Item(20, 4).Value = 1.064922260738357
$ws.Cells.Item(20, 5).Value = 1.070883222547021
$ws.Cells.Item(20, 6).Value = 1.077410647119308
$ws.Cells.Item(20, 9).Value = 1.054033970461048
$ws.Cells.Item(20, 10).Value = 1.06347996057548
$ws.Cells.Item(20, 11).Value = 1.068197773966788
$ws.Cells.Item(20, 12).Value = 1.074139253659175
$ws.Cells.Item(20, 13).Value = 1.080645615508411
$ws.Cells.Item(20, 14).Value = 1.064990224102353

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.056283167839572
$ws.Cells.Item(21, 4).Value = 1.063980994218394
$ws.Cells.Item(21, 5).Value = 1.069784498121802
$ws.Cells.Item(21, 6).Value = 1.07628535505936
$ws.Cells.Item(21, 9).Value = 1.053658180462693
$ws.Cells.Item(21, 10).Value = 1.062726205276909
$ws.Cells.Item(21, 11).Value = 1.067460018731884
$ws.Cells.Item(21, 12).Value = 1.073243280396678
$ws.Cells.Item(21, 13).Value = 1.079721747940335
$ws.Cells.Item(21, 14).Value = 1.064235398384802

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.055565858020399
$ws.Cells.Item(22, 4).Value = 1.063389528818651
$ws.Cells.Item(22, 5).Value = 1.069094472660991
$ws.Cells.Item(22, 6).Value = 1.075578659892344
$ws.Cells.Item(22, 9).Value = 1.053420729707932
$ws.Cells.Item(22, 10).Value = 1.062252047695045
$ws.Cells.Item(22, 11).Value = 1.066995809254553
$ws.Cells.Item(22, 12).Value = 1.072680065544719
$ws.Cells.Item(22, 13).Value = 1.079141033478209
$ws.Cells.Item(22, 14).Value = 1.063760567444785

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.055946077000973
$ws.Cells.Item(23, 4).Value = 1.063703042707916
$ws.Cells.Item(23, 5).Value = 1.069460193418084
$ws.Cells.Item(23, 6).Value = 1.075953214306122
$ws.Cells.Item(23, 9).Value = 1.053546718270778
$ws.Cells.Item(23, 10).Value = 1.062503430327909
$ws.Cells.Item(23, 11).Value = 1.067241928784905
$ws.Cells.Item(23, 12).Value = 1.072978624882654
$ws.Cells.Item(23, 13).Value = 1.079448866014607
$ws.Cells.Item(23, 14).Value = 1.064012307069822

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.057443860345067
$ws.Cells.Item(24, 4).Value = 1.064938047001446
$ws.Cells.Item(24, 5).Value = 1.070901656075638
$ws.Cells.Item(24, 6).Value = 1.077429526639733
$ws.Cells.Item(24, 9).Value = 1.054040250659434
$ws.Cells.Item(24, 10).Value = 1.063492593245484
$ws.Cells.Item(24, 11).Value = 1.068210136494698
$ws.Cells.Item(24, 12).Value = 1.074154276740242
$ws.Cells.Item(24, 13).Value = 1.08066110689355
$ws.Cells.Item(24, 14).Value = 1.065002874712198

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.059183797082044
$ws.Cells.Item(25, 4).Value = 1.066372678867499
$ws.Cells.Item(25, 5).Value = 1.07257778911539
$ws.Cells.Item(25, 6).Value = 1.079146249928127
$ws.Cells.Item(25, 9).Value = 1.054607832958299
$ws.Cells.Item(25, 10).Value = 1.06463939079545
$ws.Cells.Item(25, 11).Value = 1.069332128124579
$ws.Cells.Item(25, 12).Value = 1.075519051541437
$ws.Cells.Item(25, 13).Value = 1.082068508217472
$ws.Cells.Item(25, 14).Value = 1.066151300846224

